$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Source table")

# Row 54: was "Observation station file (old)" (unsupported) -> now the .xyn variant, marked supported.
$ws.Range("A54").Value = "Observation station file (.xyn)"

# Row 55: was "Observation station file (new)" -> renamed to the .ini variant; other columns unchanged.
$ws.Range("A55").Value = "Observation station file (.ini)"

$ws.Range("B54").Value = "X"
$ws.Range("C54").Value = "X"
$ws.Range("D54").Value = "0.5.0"
$ws.Range("E54").Value = "hydrolib.core.dflowm.xyn.models"
$ws.Range("F54").Value = "XYNModel"

# Restore the selection recorded for this sheet after the edit.
$ws.Activate()
$ws.Range("G56").Select()
